$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pepID")

$newIds = @(
    "PEP_ID-2009386",
    "PEP_ID-2009388",
    "PEP_ID-2009389",
    "PEP_ID-2009391",
    "PEP_ID-2009392",
    "PEP_ID-2009393"
)

$startRow = 29
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}
